$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 12.0531
$ws.Range("B8").Value = 4.7767
$ws.Range("B10").Value = 8.510600000000005
$ws.Range("B12").Value = 6.001899999999997
$ws.Range("C13").Value = -12.56889999999999
$ws.Range("B18").Value = 5.371900000000004
$ws.Range("E20").Value = 13.32449999999999
